# "Actualizar" refresh: shift the three stacked 14-row availability-check
# snapshots (rows 2-15, 16-29, 30-43) down by one slot and stamp the top
# slot with the new check timestamp. The oldest snapshot (old rows 30-43)
# falls off the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Value that was in the middle slot (D16:D29) becomes the value for the
# oldest slot (D30:D43); the newest slot (D2:D15) and the middle slot
# (D16:D29) get stamped with this run's freshly captured timestamps
# (the middle slot's new stamp carries the same sub-millisecond jitter
# the source logger produced when it re-serialised the previous "now"
# value).
$previousMiddle = $ws.Range("D16").Value2

$newTimestamp = 44233.53308902856
$shiftedTimestamp = 44233.51190966435

$ws.Range("D30:D43").Value2 = $previousMiddle
$ws.Range("D16:D29").Value2 = $shiftedTimestamp
$ws.Range("D2:D15").Value2 = $newTimestamp
